# Update "想去人数" (want-to-go count) values in column F across sheets,
# reflecting newly generated counts (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 266
$ws1.Range("F6").Value = 34
$ws1.Range("F7").Value = 262
$ws1.Range("F8").Value = 210
$ws1.Range("F9").Value = 1956
$ws1.Range("F11").Value = 4595
$ws1.Range("F12").Value = 74
$ws1.Range("F13").Value = 323

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 7

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 266
$ws4.Range("F8").Value = 34
$ws4.Range("F9").Value = 262
$ws4.Range("F10").Value = 210
$ws4.Range("F11").Value = 7
$ws4.Range("F13").Value = 1956
$ws4.Range("F15").Value = 4595
$ws4.Range("F16").Value = 74
$ws4.Range("F17").Value = 323
